# Long straddle report - add the 30-Aug-2024 NIFTY BANK dry-run rows
# (rows 6-8) to the existing report sheet, right after the existing
# NIFTY 50 rows (4-5). "fry run on real market on 30 auh for SEP
# monthly expiry"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("NIFTY BANK", "2024-08-30", 15, 1, "09:20:43", 51401.8, 876.2, 630.15, 1506.35, 22595.25, 51400, "09:20:46", 51394, 870.2, 633.4, 1503.6, 22554, -2.75, -41.25, 91.82704799999999, -133.077048, -0.5889602814750887, "MAX_LOSS_WHILE_LOOKING_FOR_INITIAL_SL"),
    @("NIFTY BANK", "2024-08-30", 15, 1, "09:24:08", 51398.3, 879.7, 626.8, 1506.5, 22597.5, 51400, "09:24:10", 51402.1, 876.4, 626.8, 1503.2, 22548, -3.300000000000182, -49.50000000000364, 91.76662279999996, -141.2666228000036, -0.6251427051665166, "MAX_LOSS_WHILE_LOOKING_FOR_INITIAL_SL"),
    @("NIFTY BANK", "2024-08-30", 15, 1, "09:35:01", 51333.35, 919, 593.45, 1512.45, 22686.75, 51300, "12:00:00", 51392.55, 946.05, 574.95, 1521, 22815, 8.549999999999955, 128.25, 93.06851119999999, 35.18148880000001, 0.1550750495333179, "TIME_ELAPSED_WHILE_LOOKING_FOR_INITIAL_SL")
)

$lastCol = 23
$startRow = 6

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]

    # Column B holds dates formatted as plain text (e.g. "2024-08-30").
    # Excel's COM layer auto-coerces ISO-date-looking strings into date
    # serials, so force text formatting on that cell before assigning it.
    $ws.Cells.Item($r, 2).NumberFormat = "@"

    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }

    # A few columns (C, D, E, L) inherit a column-level number format from
    # the sheet's <cols> definition, and column B now carries the forced
    # text format from above. Re-level the whole row back to the sheet's
    # plain default formatting (matching its unstyled neighbour cells) by
    # copying the untouched format of column A onto the row - this only
    # copies formatting, so the values just written are left alone.
    $ws.Cells.Item($r, 1).Copy()
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, $lastCol)).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
